# Updates for May 2025 application
# Adds two new "project" entries for ASR consulting work (California First 5
# County Commissions / County Offices of Education) near the bottom of the
# CV entries table, shifting the existing rows 42-50 down to 44-52.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert two fresh rows above the current row 42 -----------------------
$ws.Rows.Item(42).Insert()
$ws.Rows.Item(42).Insert()

# --- Row 42: Lead PM for six CA First 5 County Commissions -----------------
$ws.Range("A42").Value = "project"
$ws.Range("B42").Value = 2024
$ws.Range("C42").Value = "current"
$ws.Range("D42").Value = 'Lake, Kern, Siskiyou, Sutter, Yuba, and Ventura Counties, CA'
$ws.Range("E42").Value = "Applied Survey Research"
$ws.Range("F42").Value = 'Lead Project Manager for ASR''s consulting work with six California First 5 County Commissions. Support Commission efforts related to annual program evaluation, including data collection, management, analysis, and reporting. Facilitate strategic planning initiatives, including the design and distribution of surveys via Qualtrics and the collection of qualitative data through focus groups and key informant interviews. Spearhead special projects, including the design, implementation, and analysis of key surveys such as the First 5 Kern Developmental Screening Gap Assessment and the Parent/Guardian Survey and Vaccination and Immunization for the Kery County Immunization Coalition.'
$ws.Range("G42").Value = 'https://jim-asr.shinyapps.io/KernCountyImmunizationCoalition_SurveyResults/'
$ws.Rows.Item(42).RowHeight = 409.6

# --- Row 43: Lead PM for three CA County Offices of Education --------------
$ws.Range("A43").Value = "project"
$ws.Range("B43").Value = 2024
$ws.Range("C43").Value = "current"
$ws.Range("D43").Value = 'Butte, Plumas, and Santa Cruz County Offices of Education, CA'
$ws.Range("E43").Value = "Applied Survey Research"
$ws.Range("F43").Value = 'Serve as Lead Project Manager for ASR''s consulting work with three California County Offices of Education. Oversee the design and execution of Needs Assessments and parent surveys focused on perceptions, barriers, and priorities related to child care availability, priorities, and barriers to access.'
$ws.Range("G43").Value = 'https://jim-asr.shinyapps.io/SantaCruz_UPK_Survey/'
$ws.Rows.Item(43).RowHeight = 238

# --- Minor row-6 re-wrap height tweak (unrelated content, cosmetic only) ---
$ws.Rows.Item(6).RowHeight = 409.5

# --- Sheet view: scroll to the newly added rows & update selection ---------
$ws.Application.ActiveWindow.ScrollRow = 42
$ws.Range("G42").Select()

Write-Output "done"
